# Handles float input without breaking stuff
#
# Rebuilds the marksheet's score summary (rows 10-12) and the per-question
# "Student Ans" / "Correct Ans" grid (rows 15-40) now that the student's
# (float-safe) answers have been (re)graded, and drops the unused 3rd
# question block (columns G:H) plus the 2nd block's now-empty rows (19-40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Score summary (rows 10-12): Right / Wrong / Not-Attempt / Max, the
#    per-question marking weights, and the computed total.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "No."
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = 28

$ws.Range("A11").Value = "Marking"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = 44
$ws.Range("C12").Value = -5
$ws.Range("E12").Value = "39/112"

# A10:A12 get the same bordered/bold "title" look already used by A9
# (style index 4 in the saved workbook) instead of the default style.
$ws.Range("A9").Copy()
foreach ($addr in @("A10", "A11", "A12")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 2. Per-question grid (rows 16-40): fill in the student's answers in
#    column A (and column D for question 2), recoloring each cell green
#    ("correctStyle", already used as style 5) when it matches the
#    correct answer, red ("incorrectStyle", style 6) when it doesn't, and
#    leaving the existing empty style (7) for not-attempted questions.
# ---------------------------------------------------------------------
$correct = @{
    "A16" = "Option C"; "A18" = "Option B"; "A19" = "Option C";
    "A22" = "Option D"; "A24" = "Option A"; "A26" = "Option C";
    "A27" = "Option C"; "A28" = "Option D"; "A29" = "Option D";
    "A33" = "Option D"; "A34" = "Option A"; "A35" = "Option D";
    "A38" = "Option A";
    "D16" = "Option A"; "D17" = "Option B"; "D18" = "Option B";
}
foreach ($addr in $correct.Keys) {
    $ws.Range($addr).Value = $correct[$addr]
}

# Cells whose student answer matches the correct answer -> green (style 5,
# sourced from B10 which already carries that style).
$ws.Range("B10").Copy()
foreach ($addr in @("A18", "A19", "A22", "A24", "A26", "A28", "A29", "A33", "A35", "A38", "D16")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

# Cells whose student answer is wrong -> red (style 6, sourced from C10).
$ws.Range("C10").Copy()
foreach ($addr in @("A16", "A27", "A34", "D17", "D18")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 3. Drop the unused 3rd question block entirely, and the 2nd question
#    block's rows beyond its last real question.
# ---------------------------------------------------------------------
$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

Write-Host "done"
